$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark additional checklist items as "done" (green) by re-using the
#     same fill already applied to other completed rows, via a
#     copy/paste-special-formats so the workbook's existing theme fill
#     is reused instead of a brand new raw RGB fill being created. ---

# A11 ("mejla ut en länk ...") -> same look as A4/A7 (green, no vertical-center)
$ws.Range("A4").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# A13 ("statistik över testresultatet ...") -> same look as A2/A3 (green, centered)
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# A21 / A22 ("När man skiftar ...", "fel medd ...") -> same look as A18/A19 (green)
$ws.Range("A18").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- New checklist item: "Dela till gruppen" (Share to the group) tied to row 12 ---
$ws.Range("D12").Value = "Dela till gruppen"

# --- Restore the last-used selection as recorded by the editing session ---
$ws.Range("A26:A27").Select()
